# Insert a new data row at row 25 (shifting all subsequent rows down by one,
# which grows the sheet from 141 to 142 data rows / A1:R141 to A1:R142),
# then populate the newly inserted row with its values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 25..141 down to 26..142
$ws.Rows("25:25").Insert()

# Populate the new row 25 with the inserted record's data
$ws.Range("A25").Value = 10
$ws.Range("B25").Value = "Vega Modelo de Temuco"
$ws.Range("C25").Value = "La Araucanía"
$ws.Range("D25").Value = 45250
$ws.Range("E25").Value = 9
$ws.Range("F25").Value = 300000001
$ws.Range("G25").Value = "Rabanito"
$ws.Range("H25").Value = "Sin especificar"
$ws.Range("I25").Value = "Primera"
$ws.Range("J25").Value = 70
$ws.Range("K25").Value = 8000
$ws.Range("L25").Value = 9000
$ws.Range("M25").Value = 8286
$ws.Range("N25").Value = "$/docena de paquetes"
$ws.Range("O25").Value = "Provincia de Cautín"
$ws.Range("P25").Value = 690
$ws.Range("Q25").Value = 12
$ws.Range("R25").Value = "Hortaliza"
